$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 441.91666
$ws.Range("J2").Value = 716.6667
$ws.Range("L2").Value = 716.6667
$ws.Range("N2").Value = -942.6667
$ws.Range("H15").Value = 606005.25
$ws.Range("I15").Value = 606005.25
$ws.Range("K15").Value = 1818015.75
$ws.Range("M15").Value = -1817846.75
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H42").Value = 104.77778
$ws.Range("I42").Value = 90
$ws.Range("J42").Value = 123.25
$ws.Range("K42").Value = 270
$ws.Range("L42").Value = 369.75
$ws.Range("M42").Value = -40
$ws.Range("N42").Value = -829.75
$ws.Range("H55").Value = 169.6
$ws.Range("I55").Value = 169.6
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 169.6
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = 44.40000000000001
$ws.Range("N55").ClearContents()
$ws.Range("H69").Value = 13959.533
$ws.Range("I69").Value = 13277.111
$ws.Range("J69").Value = 14983.167
$ws.Range("K69").Value = 39831.333
$ws.Range("L69").Value = 44949.501
$ws.Range("M69").Value = -38957.333
$ws.Range("N69").Value = -46697.501
$ws.Range("H72").Value = 13959.533
$ws.Range("I72").Value = 13277.111
$ws.Range("J72").Value = 14983.167
$ws.Range("K72").Value = 119493.999
$ws.Range("L72").Value = 134848.503
$ws.Range("M72").Value = -115125.999
$ws.Range("N72").Value = -143584.503
$ws.Range("H74").Value = 4529.6
$ws.Range("I74").Value = 4529.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4529.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3593.6
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 4529.6
$ws.Range("I77").Value = 4529.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 22648
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -17968
$ws.Range("N77").ClearContents()
$ws.Range("H138").Value = 153804.94
$ws.Range("I138").Value = 550068.3
$ws.Range("J138").Value = 5206.1562
$ws.Range("K138").Value = 1650204.9
$ws.Range("L138").Value = 15618.4686
$ws.Range("M138").Value = -1645064.9
$ws.Range("N138").Value = -25898.4686

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 76272.57000000001
$ws.Range("I45").Value = 91143.44
$ws.Range("J45").Value = 7866.6
$ws.Range("K45").Value = 91143.44
$ws.Range("L45").Value = 7866.6
$ws.Range("M45").Value = -90766.44
$ws.Range("N45").Value = -8620.6
$ws.Range("H74").Value = 3830.2632
$ws.Range("I74").Value = 12972.5
$ws.Range("K74").Value = 12972.5
$ws.Range("M74").Value = -12098.5
$ws.Range("H77").Value = 3830.2632
$ws.Range("I77").Value = 12972.5
$ws.Range("K77").Value = 64862.5
$ws.Range("M77").Value = -60494.5
$ws.Range("H122").Value = 1435154
$ws.Range("I122").Value = 6576.5625
$ws.Range("K122").Value = 19729.6875
$ws.Range("M122").Value = -17279.6875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 5907.933
$ws.Range("J64").Value = 1047
$ws.Range("L64").Value = 1047
$ws.Range("N64").Value = -1497
$ws.Range("H67").Value = 5907.933
$ws.Range("J67").Value = 1047
$ws.Range("L67").Value = 1047
$ws.Range("N67").Value = -2607

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 127037.234
$ws.Range("I105").Value = 134664.56
$ws.Range("K105").Value = 134664.56
$ws.Range("M105").Value = -132917.56
$ws.Range("H132").Value = 5881.3335
$ws.Range("I132").Value = 6216
$ws.Range("K132").Value = 18648
$ws.Range("M132").Value = -16118
$ws.Range("H134").Value = 3319.0667
$ws.Range("I134").Value = 3234.7144
$ws.Range("J134").Value = 4500
$ws.Range("K134").Value = 9704.143199999999
$ws.Range("L134").Value = 13500
$ws.Range("M134").Value = -7169.143199999999
$ws.Range("N134").Value = -18570
$ws.Range("H141").Value = 407702.88
$ws.Range("J141").Value = 461313.9
$ws.Range("L141").Value = 461313.9
$ws.Range("N141").Value = -471673.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 133.66667
$ws.Range("I12").Value = 133.66667
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 401.00001
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -228.00001
$ws.Range("N12").ClearContents()
$ws.Range("H113").Value = 1379.3077
$ws.Range("I113").Value = 911.1667
$ws.Range("J113").Value = 1780.5714
$ws.Range("K113").Value = 2733.5001
$ws.Range("L113").Value = 5341.7142
$ws.Range("M113").Value = -563.5001000000002
$ws.Range("N113").Value = -9681.7142
$ws.Range("H122").Value = 5344.6387
$ws.Range("I122").Value = 1235.25
$ws.Range("J122").Value = 7399.3335
$ws.Range("K122").Value = 11117.25
$ws.Range("L122").Value = 66594.0015
$ws.Range("M122").Value = -8667.25
$ws.Range("N122").Value = -71494.0015
$ws.Range("H140").Value = 770904.0600000001
$ws.Range("I140").Value = 770904.0600000001
$ws.Range("K140").Value = 2312712.18
$ws.Range("M140").Value = -2307532.18
$ws.Range("H141").Value = 1865.6
$ws.Range("I141").Value = 1776
$ws.Range("J141").Value = 2000
$ws.Range("K141").Value = 5328
$ws.Range("L141").Value = 6000
$ws.Range("M141").Value = -148
$ws.Range("N141").Value = -16360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 4883.5
$ws.Range("I99").Value = 4883.5
$ws.Range("K99").Value = 4883.5
$ws.Range("M99").Value = -2637.5
$ws.Range("H122").Value = 9760.267
$ws.Range("I122").Value = 6270.2173
$ws.Range("K122").Value = 18810.6519
$ws.Range("M122").Value = -16360.6519
$ws.Range("H123").Value = 47000
$ws.Range("J123").Value = 47000
$ws.Range("L123").Value = 47000
$ws.Range("N123").Value = -51900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 679.8889
$ws.Range("I22").Value = 597.5
$ws.Range("J22").Value = 745.8
$ws.Range("K22").Value = 597.5
$ws.Range("L22").Value = 745.8
$ws.Range("M22").Value = -302.5
$ws.Range("N22").Value = -1335.8
$ws.Range("H27").Value = 679.8889
$ws.Range("I27").Value = 597.5
$ws.Range("J27").Value = 745.8
$ws.Range("K27").Value = 597.5
$ws.Range("L27").Value = 745.8
$ws.Range("M27").Value = -490.5
$ws.Range("N27").Value = -959.8
$ws.Range("H55").Value = 4190.636
$ws.Range("I55").Value = 1299.8334
$ws.Range("J55").Value = 7659.6
$ws.Range("K55").Value = 1299.8334
$ws.Range("L55").Value = 7659.6
$ws.Range("M55").Value = -1126.8334
$ws.Range("N55").Value = -8005.6
$ws.Range("H122").Value = 9200
$ws.Range("I122").Value = 13650
$ws.Range("J122").Value = 7222.222
$ws.Range("K122").Value = 40950
$ws.Range("L122").Value = 21666.666
$ws.Range("M122").Value = -38500
$ws.Range("N122").Value = -26566.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 76927590
$ws.Range("H4").Value = 1344.375
$ws.Range("J4").Value = 1627.6923
$ws.Range("L4").Value = 1627.6923
$ws.Range("N4").Value = -1853.6923
$ws.Range("H5").Value = 12515000
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H26").Value = 9491.5
$ws.Range("I26").Value = 9491.5
$ws.Range("K26").Value = 9491.5
$ws.Range("M26").Value = -9198.5
$ws.Range("H40").Value = 52995
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 52995
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 52995
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -53293
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("H58").Value = 6811334
$ws.Range("I58").Value = 14167.5
$ws.Range("K58").Value = 14167.5
$ws.Range("M58").Value = -13859.5
$ws.Range("H135").Value = 64177.5
$ws.Range("J135").Value = 62236.668
$ws.Range("L135").Value = 62236.668
$ws.Range("N135").Value = -72376.66800000001
